$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("phen_oncox")

# Disease Ontology (row 3) source_version: v2024-09-27 -> v2024-11-01
$ws.Range("E3").Value = "v2024-11-01"

# Experimental Factor Ontology (row 4) source_version: v3.71.0 -> v3.72.0
$ws.Range("E4").Value = "v3.72.0"
